$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 429, shifting existing rows 429:545 down to 430:546
$ws.Rows("429:429").Insert()

# Populate the newly inserted row 429 with the new price record
$ws.Range("A429").Value = 11
$ws.Range("B429").Value = "Vega Monumental Concepción"
$ws.Range("C429").Value = "Bíobío"
$ws.Range("D429").Value = 44785
$ws.Range("E429").Value = 8
$ws.Range("F429").Value = "Fruta"
$ws.Range("G429").Value = 100108
$ws.Range("H429").Value = "Tropicales y subtropicales"
$ws.Range("I429").Value = 100108006
$ws.Range("J429").Value = "Plátano"
$ws.Range("K429").Value = "Sin especificar"
$ws.Range("L429").Value = "Pintón"
$ws.Range("M429").Value = 1050
$ws.Range("N429").Value = 19000
$ws.Range("O429").Value = 20000
$ws.Range("P429").Value = 19524
$ws.Range("Q429").Value = "`$/caja 20 kilos"
$ws.Range("R429").Value = "Ecuador"
$ws.Range("S429").Value = 976
$ws.Range("T429").Value = 20
